$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.17390000000001
$ws.Range("A21").Value = -20.06099999999998
$ws.Range("A23").Value = -20.20879999999998
$ws.Range("A25").Value = -21.84089999999999
$ws.Range("A53").Value = -21.9389
$ws.Range("A57").Value = -22.61160000000002
$ws.Range("A59").Value = -22.13019999999999
$ws.Range("A69").Value = -21.54919999999998
$ws.Range("A79").Value = -20.1764
$ws.Range("A83").Value = -21.82830000000001
$ws.Range("A93").Value = -21.44090000000001
